$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout: A date | B fly_time | C fuel consumption litres | D fuel consumption cost kDKK | E last activity
# New layout:      A date | B fly_time | C fly_cost | D fuel consumption litres | E fuel consumption cost kDKK | F day cost per all people | G last activity

# --- 1) Insert the two new columns (fly_cost, day cost per all people) ---
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("F1").EntireColumn.Insert()

# --- 2) Insert 3 new rows above the existing "total" row (currently row 7) for the new
#         weather-delay days; this shifts the old total row (and its cells, including the
#         blank trailing "last activity" cell) down to row 10 ---
$ws.Range("A7:A9").EntireRow.Insert()

# --- 3) Insert 2 more rows after the total row for the two new grand-total summary lines ---
$ws.Range("A11:A12").EntireRow.Insert()

# --- Header row ---
$ws.Range("A1").Value = "date (YYYY-MM-DD)"
$ws.Range("B1").Value = "fly_time"
$ws.Range("C1").Value = "fly_cost"
$ws.Range("D1").Value = "fuel consumption litres"
$ws.Range("E1").Value = "fuel consumption cost kDKK"
$ws.Range("F1").Value = "day cost per all people"
$ws.Range("G1").Value = "last activity"

# Apply the bordered/bold/centered header style to the two newly inserted header cells
# (copy formats only from an existing header cell so the new headers match visually)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Daily rows: 5 original flying days (now with fly_cost + day cost added) plus 3 new weather-delay days ---

$rows = @(
    @{ r = 2; date = "2021-08-20"; fly_time = 3.3; fly_cost = 44.3; fuel_l = 2478; fuel_c = 25.1;               day_cost = 15;  act = "move people and cargo to QAN airport" },
    @{ r = 3; date = "2021-08-21"; fly_time = 2.9; fly_cost = 39.1; fuel_l = 1754; fuel_c = 17.8;               day_cost = 15;  act = "return to QAN to overnight" },
    @{ r = 4; date = "2021-08-22"; fly_time = 3.6; fly_cost = 48.4; fuel_l = 3097; fuel_c = 31.4;               day_cost = 15;  act = "return to QAAN to overnight" },
    @{ r = 5; date = "2021-08-23"; fly_time = 1.6; fly_cost = 21.5; fuel_l = 964;  fuel_c = 9.800000000000001; day_cost = 15;  act = "return to QAN to overnight" },
    @{ r = 6; date = "2021-08-24"; fly_time = 3.3; fly_cost = 44.3; fuel_l = 1488; fuel_c = 15.1;               day_cost = 15;  act = "return to YRB with no PAX" },
    @{ r = 7; date = "2021-08-26"; fly_time = 4;   fly_cost = 53.6; fuel_l = 0;    fuel_c = 0;                  day_cost = 7.5; act = "wx delay day" },
    @{ r = 8; date = "2021-08-27"; fly_time = 4;   fly_cost = 53.6; fuel_l = 0;    fuel_c = 0;                  day_cost = 7.5; act = "wx delay day" },
    @{ r = 9; date = "2021-08-28"; fly_time = 4;   fly_cost = 53.6; fuel_l = 0;    fuel_c = 0;                  day_cost = 7.5; act = "wx delay day" }
)

foreach ($row in $rows) {
    $r = $row.r

    # Keep the date column as plain text (e.g. "2021-08-20"), not an auto-converted date serial.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row.date

    $ws.Cells.Item($r, 2).Value = $row.fly_time
    $ws.Cells.Item($r, 3).Value = $row.fly_cost
    $ws.Cells.Item($r, 4).Value = $row.fuel_l
    $ws.Cells.Item($r, 5).Value = $row.fuel_c
    $ws.Cells.Item($r, 6).Value = $row.day_cost
    $ws.Cells.Item($r, 7).Value = $row.act
}

# --- Totals row (row 10 after the inserts above); note G10 (old E7) is left untouched,
#     preserving its original blank "last activity" cell. ---
$ws.Range("A10").Value = "total"
$ws.Range("B10").Value = 26.7
$ws.Range("C10").Value = 358.4
$ws.Range("D10").Value = 9781
$ws.Range("E10").Value = 99.2
$ws.Range("F10").Value = 97.5

# --- Grand total summary rows ---
$ws.Range("A11").Value = "grand total (MDKK)"
$ws.Range("B11").Value = 0.5551

$ws.Range("A12").Value = "grand total incl. quarantine (MDKK)"
$ws.Range("B12").Value = 0.6301
